# Add a new "date" column (D) with formatted date/time strings for each
# data row, matching the "new design implement and some functionality
# changes" commit: a payment table gains a date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("D1").Value = "date"

# Row 2 has a date value containing an embedded line break, so it needs
# wrap text on and a taller row to show both lines.
$ws.Range("D2").Value = "01 mar 2023, 06:09 pm`n"
$ws.Range("D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 24.3

# Rows 3 and 4 are plain single-line date strings.
$ws.Range("D3").Value = "08 mar 2023, 02:19 am"
$ws.Range("D4").Value = "10 mar 2023, 12:00 pm"

# Leave the final selection on the last cell touched, as in the edit.
$null = $ws.Range("D4").Select()
